$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set - words broken out individually with updated coordinates/sizes.
# Columns: A=text, B=x, C=y, D=width, E=height
$data = @(
    @("Greek ", 152, 644, 58, 23),
    @("mythology ", 210, 644, 93, 23),
    @("gift ", 870, 926.8, 31, 23),
    @("of ", 901, 926.8, 22, 23),
    @("prophecy. ", 152, 951.8, 89, 23),
    @("Trojan ", 356, 1001.8, 59, 23),
    @("Horse ", 415, 1001.8, 57, 23),
    @("trick, ", 472, 1001.8, 46, 23),
    @("Agamemnon" + [char]0x2019 + "s ", 226, 1026.8, 128, 23),
    @("Bronze ", 349, 1193.2, 67, 23),
    @("Age. ", 416, 1193.2, 45, 23)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $row++
}
